$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The comment that currently sits on A7 ("Error handling strategy...") needs
# to end up anchored one row higher (on A6) once row 2 is removed below, since
# the runtime does not auto-shift comment anchors on row delete. Recreate it
# on the new target cell before shifting rows so its text is preserved.
$oldComment = $ws.Range("A7").Comment
$commentText = $oldComment.Text()
$oldComment.Delete()

# Delete row 2 ("Textured surfaces" / estimate 2) entirely; rows 3-14 shift up
# by one, so the old row 7 becomes row 6, matching where the comment must end
# up (A6).
$ws.Rows.Item(2).Delete()

$newComment = $ws.Range("A6").AddComment($commentText)

# Restore the selection state to match the post-edit workbook (entire row 2
# selected, as after a row deletion via the row header).
$ws.Range("A2:XFD2").Select()
